# Validating the GN date range in Excel.
# The "DATE" column on "Mismatched Data" and the "DATE" column on
# "not_in_Portal" were stored as plain text strings (e.g. "2025-02-14").
# Convert them into real Excel date values formatted as YYYY-MM-DD so the
# workbook can do real date-range validation/comparisons on them.

$wb = $excel.ActiveWorkbook

# --- "Mismatched Data" sheet: single data row, DATE is column S (row 2) ---
$wsMismatched = $wb.Worksheets.Item("Mismatched Data")
$wsMismatched.Range("S2").NumberFormat = "YYYY-MM-DD"
$wsMismatched.Range("S2").Value = "2025-02-14"

# --- "not_in_Portal" sheet: DATE is column E, rows 2-101 ---
$wsNotInPortal = $wb.Worksheets.Item("not_in_Portal")

$wsNotInPortal.Range("E2:E42").NumberFormat = "YYYY-MM-DD"
$wsNotInPortal.Range("E2:E42").Value = "2025-02-19"

$wsNotInPortal.Range("E43:E101").NumberFormat = "YYYY-MM-DD"
$wsNotInPortal.Range("E43:E101").Value = "2025-02-18"
